{"js": "// Locate the run of text that ends the \"Add multiple bookings...\" paragraph\n// and append the new sentence fragment right after it (this keeps the\n// insertion before the trailing _GoBack bookmark, matching the target diff).\nconst body = context.document.body;\nconst found = body.search(\n  \"if there are any limits/ issues when multiple bookings are recorded.\",\n  { matchCase: true }\n);\nfound.load(\"text\");\nawait context.sync();\n\nif (found.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to extend.\");\n}\n\nconst tail = found.items[0];\ntail.insertText(\" (via selenium or jscript in console)\", Word.InsertLocation.after);\nawait context.sync();\n\n// Re-resolve the owning paragraph (content changed, so re-query by text).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Add multiple bookings and see\") !== -1) {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not re-locate the target paragraph.\");\n}\n\n// Insert a new, empty paragraph right after it \u2026\nconst blankParagraph = target.insertParagraph(\"\", Word.InsertLocation.after);\n\n// \u2026 followed by a paragraph containing the new note.\nblankParagraph.insertParagraph(\"Run Pa11y for accessible\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Find the sentence that ends the \"Add multiple bookings...\" paragraph and\n# append the new fragment right after it (lands before the trailing\n# _GoBack bookmark, matching the target diff).\n$rng = $d.Content\n$found = $rng.Find.Execute(\"if there are any limits/ issues when multiple bookings are recorded.\")\n$rng.Collapse(0)\n$rng.InsertAfter(\" (via selenium or jscript in console)\")\n\n# Re-locate the owning paragraph now that its text has changed.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Add multiple bookings and see*\") {\n        $target = $p\n        break\n    }\n}\n\n# Insert a new, empty paragraph right after it \u2026\n$target.Range.InsertParagraphAfter()\n\n# \u2026 followed by a paragraph containing the new note.\n$blankParagraph = $target.Next()\n$blankParagraph.Range.InsertParagraphAfter()\n$newParagraph = $blankParagraph.Next()\n$newParagraph.Range.Text = \"Run Pa11y for accessible\"\n"}
